$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Proposta de Desnormalizacao" (K) and "Analise de Volume de
# Transacoes" (J) columns as done ("ok") on row 5, same as the existing
# F5:I5 marks.
$ws.Range("J5").Value = "ok"
$ws.Range("K5").Value = "ok"

# Move the active selection/view to the newly-filled K5 cell and scroll
# the window right so the edited columns are in view.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K5").Select()
